# Generate Report for Handoff
#
# The d74854b4-26bb-4c7f-a523-8082fbbe6f40.md file has moved from
# "Handed back" to "Ready for handoff" status, with updated handoff
# timestamps, and the zh-cn / de-de locale sheets now record that the
# handback file that was produced is stale (an Error Detail message),
# plus the Error Detail column is widened so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c24e5a413eb6017ef38dccc6c1f8676bc1cbad65/e2e/d74854b4-26bb-4c7f-a523-8082fbbe6f40.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/813f9171d71aadbf2e4037f34e45132490c6e0b7/e2e/d74854b4-26bb-4c7f-a523-8082fbbe6f40.md."

# --- Overview sheet: row for d74854b4-26bb-4c7f-a523-8082fbbe6f40.md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-31 18:55:50"

# --- zh-cn sheet: row for d74854b4-26bb-4c7f-a523-8082fbbe6f40.md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-31 18:55:45"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").ColumnWidth = 39.17

# --- de-de sheet: row for d74854b4-26bb-4c7f-a523-8082fbbe6f40.md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-31 18:55:50"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.17
